$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.923.07'
Set-TextValue $ws.Range("E2") '  +0.04%  '

Set-TextValue $ws.Range("D3") '1.549.58'
Set-TextValue $ws.Range("E3") '  +0.27%  '

Set-TextValue $ws.Range("E4") '  -0.18%  '

Set-TextValue $ws.Range("D5") '206.63'
Set-TextValue $ws.Range("E5") '  +0.45%  '

Set-TextValue $ws.Range("D6") '0.487'
Set-TextValue $ws.Range("E6") '  +0.84%  '

Set-TextValue $ws.Range("E7") '  -0.22%  '

Set-TextValue $ws.Range("D8") '22.06'
Set-TextValue $ws.Range("E8") '  +3.10%  '

Set-TextValue $ws.Range("E9") '  -0.40%  '

Set-TextValue $ws.Range("E11") '  -0.26%  '

Set-TextValue $ws.Range("D12") '1.771.13'
Set-TextValue $ws.Range("E12") '  +0.28%  '

Set-TextValue $ws.Range("D13") '1.554.83'
Set-TextValue $ws.Range("E13") '  +0.65%  '

Set-TextValue $ws.Range("D14") '3.73'
Set-TextValue $ws.Range("E14") '  +1.25%  '

Set-TextValue $ws.Range("E15") '  +1.47%  '

Set-TextValue $ws.Range("D16") '26.923.98'
Set-TextValue $ws.Range("E16") '  +0.10%  '

Set-TextValue $ws.Range("D17") '61.59'
Set-TextValue $ws.Range("E17") '  +0.07%  '

Set-TextValue $ws.Range("D18") '217.30'
Set-TextValue $ws.Range("E18") '  +1.84%  '

Set-TextValue $ws.Range("D19") '0.0₃0695'
Set-TextValue $ws.Range("E19") '  +1.73%  '

Set-TextValue $ws.Range("E20") '  +1.16%  '

Set-TextValue $ws.Range("E21") '  -0.19%  '

Set-TextValue $ws.Range("E22") '  +0.64%  '

Set-TextValue $ws.Range("D23") '9.18'
Set-TextValue $ws.Range("E23") '  +0.15%  '

Set-TextValue $ws.Range("D24") '1.95'
Set-TextValue $ws.Range("E24") '  +0.57%  '

Set-TextValue $ws.Range("D25") '154.28'
Set-TextValue $ws.Range("E25") '  +0.77%  '

Set-TextValue $ws.Range("E26") '  -0.20%  '

Set-TextValue $ws.Range("D27") '14.91'
Set-TextValue $ws.Range("E27") '  +0.63%  '

Set-TextValue $ws.Range("E28") '  +0.96%  '

Set-TextValue $ws.Range("E29") '  -0.14%  '

Set-TextValue $ws.Range("E30") '  +1.85%  '

Set-TextValue $ws.Range("D31") '1.08'
Set-TextValue $ws.Range("E31") '  -0.78%  '

Set-TextValue $ws.Range("E32") '  -0.14%  '

Set-TextValue $ws.Range("D33") '1.418.67'
Set-TextValue $ws.Range("E33") '  +3.83%  '

Set-TextValue $ws.Range("D34") '3.06'
Set-TextValue $ws.Range("E34") '  +3.98%  '

Set-TextValue $ws.Range("E35") '  +3.16%  '

Set-TextValue $ws.Range("D36") '0.968'
Set-TextValue $ws.Range("E36") '  -0.13%  '

Set-TextValue $ws.Range("E37") '  +0.14%  '

Set-TextValue $ws.Range("E38") '  +0.41%  '

Set-TextValue $ws.Range("E39") '  +0.80%  '

Set-TextValue $ws.Range("B40") 'ARBITRUM'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D40") '0.807'
Set-TextValue $ws.Range("E40") '  +0.36%  '

Set-TextValue $ws.Range("B41") 'FraxShare'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D41") '5.76'
Set-TextValue $ws.Range("E41") '  +4.95%  '

Set-TextValue $ws.Range("E42") '  -0.19%  '

Set-TextValue $ws.Range("D43") '2.32'
Set-TextValue $ws.Range("E43") '  +4.29%  '

Set-TextValue $ws.Range("E44") '  +0.55%  '

Set-TextValue $ws.Range("D45") '64.26'
Set-TextValue $ws.Range("E45") '  +1.41%  '

Set-TextValue $ws.Range("E46") '  +0.74%  '

Set-TextValue $ws.Range("D47") '1.684.78'
Set-TextValue $ws.Range("E47") '  +0.28%  '

Set-TextValue $ws.Range("D48") '87.60'
Set-TextValue $ws.Range("E48") '  +1.53%  '

Set-TextValue $ws.Range("D49") '0.0518'
Set-TextValue $ws.Range("E49") '  +2.51%  '

Set-TextValue $ws.Range("E50") '  +4.72%  '

Set-TextValue $ws.Range("D51") '0.0952'
Set-TextValue $ws.Range("E51") '  +0.29%  '

